$d = $word.ActiveDocument

# Colors the next occurrence of $searchText at/after $fromPos red, returns the
# position right after the match (so subsequent searches continue forward).
function Color-NextMatch($searchText, $fromPos) {
    $r = $d.Range($fromPos, $d.Content.End)
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Font.Color = 255
        return $r.End
    }
    return $fromPos
}

$pos = 0

# Paragraph 1: "Farshid is very happy And satisfied." -> "Farshid " turns red.
$pos = Color-NextMatch "Farshid " $pos

# Paragraph 2: several words/phrases turn red, left to right through the text.
$pos = Color-NextMatch "weather " $pos
$pos = Color-NextMatch "worlds energy " $pos
$pos = Color-NextMatch "wonder " $pos
$pos = Color-NextMatch "much " $pos
$pos = Color-NextMatch "nature " $pos
$pos = Color-NextMatch "world" $pos
$pos = Color-NextMatch "life" $pos
$pos = Color-NextMatch "every day. " $pos
$pos = Color-NextMatch "girl" $pos
$pos = Color-NextMatch "s" $pos
$pos = Color-NextMatch "girls " $pos

# The document's "_GoBack" bookmark originally sat mid-word, inside "going"
# (splitting it into "goi" / "ng and "). After the text above is recolored,
# Word leaves that bookmark at the location of the most recent edit, right
# after "...attract girls" rather than inside "going". Recreate that here:
# first force a genuine text edit across the old split point so the run
# re-merges into a single "". I am very calm and easy going and " run, then
# drop the bookmark back in at its new location.
$mergeStart = $d.Paragraphs(2).Range.Start + 255
$mergeEnd = $d.Paragraphs(2).Range.Start + 291
$mergeRange = $d.Range($mergeStart, $mergeEnd)
$originalText = $mergeRange.Text
$mergeRange.Text = $originalText + "x"
$growRange = $d.Range($mergeStart, $mergeEnd + 1)
$growRange.Text = $originalText

$bookmarkRange = $d.Range($mergeStart, $mergeStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
